$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 289) holds a date serial value that was updated
# from 45189 (2023-09-20) to 45190 (2023-09-21) for every data row.
$ws.Range("C2:C289").Value = 45190
